## Weekly data refresh: a new price record for "Albahaca" (Vega Modelo de
## Temuco) was inserted as the new row 128, pushing the existing rows
## 128-149 down to 129-150 (dimension grows from A1:R149 to A1:R150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 128 - this shifts rows 128:149 down to 129:150 and
# carries formatting (incl. the date number-format on column D) down from
# the row above, matching how the existing rows are styled.
$ws.Rows(128).Insert()

# Populate the newly inserted row 128 with the new week's record.
$ws.Range("A128").Value = 10
$ws.Range("B128").Value = "Vega Modelo de Temuco"
$ws.Range("C128").Value = "La Araucanía"
$ws.Range("D128").Value = 44522
$ws.Range("E128").Value = 9
$ws.Range("F128").Value = 100112052
$ws.Range("G128").Value = "Albahaca"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 40
$ws.Range("K128").Value = 3500
$ws.Range("L128").Value = 4000
$ws.Range("M128").Value = 3750
$ws.Range("N128").Value = '$/paquete'
$ws.Range("O128").Value = "Región del Maule"
$ws.Range("P128").Value = 3750
$ws.Range("Q128").Value = 1
$ws.Range("R128").Value = "Hortaliza"
